$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 40; existing rows 40-80 shift down to 41-81.
$ws.Rows.Item(40).Insert()

# Populate the newly inserted row 40 with the new record's data.
$ws.Cells.Item(40, 1).Value = 10
$ws.Cells.Item(40, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(40, 3).Value = "La Araucanía"
$ws.Cells.Item(40, 4).Value = 44484
$ws.Range("D40").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(40, 5).Value = 9
$ws.Cells.Item(40, 6).Value = 100112031
$ws.Cells.Item(40, 7).Value = "Poroto verde"
$ws.Cells.Item(40, 8).Value = "Sin especificar"
$ws.Cells.Item(40, 9).Value = "Primera"
$ws.Cells.Item(40, 10).Value = 50
$ws.Cells.Item(40, 11).Value = 3000
$ws.Cells.Item(40, 12).Value = 3000
$ws.Cells.Item(40, 13).Value = 3000
$ws.Cells.Item(40, 14).Value = "$/kilo"
$ws.Cells.Item(40, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(40, 16).Value = 3000
$ws.Cells.Item(40, 17).Value = 1
$ws.Cells.Item(40, 18).Value = "Hortaliza"
